$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column L (12th column). This shifts all
# columns from L onward one position to the right, carrying along their
# cell values/styles (matches the "Тўлов тури" payment-type column being
# added right after the existing "Пробирка рақами" column).
$ws.Columns("L:L").Insert()

# Give the freshly inserted column its own width (it no longer should
# look like a best-fit copy of its former neighbour) and put the new
# header label in the header row (row 4).
$ws.Columns("L:L").ColumnWidth = 15.33
$ws.Range("L4").Value = "Тўлов тури"

# Update the hidden _FilterDatabase defined name so it spans through the
# new last column (X) instead of the old one (W).
foreach ($n in $wb.Names) {
    if ($n.Name -eq "TDSheet!_FilterDatabase") {
        $n.RefersTo = "=TDSheet!`$A`$4:`$X`$4"
    }
}

# Restore the selected cell that shifted from K5 to L5 because of the
# inserted column.
$ws.Range("L5").Select() | Out-Null
